$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2,7).Value = 4.592460999999999
$ws.Cells.Item(2,8).Value = 13.777383
$ws.Cells.Item(2,9).Value = 0.003302946473568516
$ws.Cells.Item(2,10).Value = 0.003302946473568516
$ws.Cells.Item(2,11).Value = 3
$ws.Cells.Item(2,12).Value = 1
$ws.Cells.Item(2,13).Value = 2.027115333333333
$ws.Cells.Item(2,14).Value = 6.081346
$ws.Cells.Item(2,15).Value = 0.006596284565418616
$ws.Cells.Item(2,16).Value = 0.006596284565418615
$ws.Cells.Item(2,17).Value = 9.309448110835332
$ws.Cells.Item(2,18).Value = 83.78503299751799
$ws.Cells.Item(2,19).Value = 0.00002178717484400384
$ws.Cells.Item(2,20).Value = 0.00002178717484400384
$ws.Cells.Item(3,7).Value = 4.592460999999999
$ws.Cells.Item(3,8).Value = 13.777383
$ws.Cells.Item(3,9).Value = 0.003302946473568516
$ws.Cells.Item(3,10).Value = 0.003302946473568516
$ws.Cells.Item(3,11).Value = 3
$ws.Cells.Item(3,12).Value = 1
$ws.Cells.Item(3,13).Value = 256.4443053333333
$ws.Cells.Item(3,14).Value = 769.332916
$ws.Cells.Item(3,15).Value = 0.8344762556643375
$ws.Cells.Item(3,16).Value = 0.8344762556643374
$ws.Cells.Item(3,17).Value = 1177.710470915425
$ws.Cells.Item(3,18).Value = 10599.39423823883
$ws.Cells.Item(3,19).Value = 0.002756230405923183
$ws.Cells.Item(3,20).Value = 0.002756230405923183
$ws.Cells.Item(4,7).Value = 4.592460999999999
$ws.Cells.Item(4,8).Value = 13.777383
$ws.Cells.Item(4,9).Value = 0.003302946473568516
$ws.Cells.Item(4,10).Value = 0.003302946473568516
$ws.Cells.Item(4,11).Value = 3
$ws.Cells.Item(4,12).Value = 1
$ws.Cells.Item(4,13).Value = 48.84026566666667
$ws.Cells.Item(4,14).Value = 146.520797
$ws.Cells.Item(4,15).Value = 0.158927459770244
$ws.Cells.Item(4,16).Value = 0.158927459770244
$ws.Cells.Item(4,17).Value = 224.2970153038057
$ws.Cells.Item(4,18).Value = 2018.673137734251
$ws.Cells.Item(4,19).Value = 0.0005249288928013296
$ws.Cells.Item(4,20).Value = 0.0005249288928013295
$ws.Cells.Item(5,7).Value = 1205.102620666667
$ws.Cells.Item(5,8).Value = 3615.307862
$ws.Cells.Item(5,9).Value = 0.8667225374846176
$ws.Cells.Item(5,10).Value = 0.8667225374846176
$ws.Cells.Item(5,11).Value = 3
$ws.Cells.Item(5,12).Value = 1
$ws.Cells.Item(5,13).Value = 2.027115333333333
$ws.Cells.Item(5,14).Value = 6.081346
$ws.Cells.Item(5,15).Value = 0.006596284565418616
$ws.Cells.Item(5,16).Value = 0.006596284565418615
$ws.Cells.Item(5,17).Value = 2442.882000593584
$ws.Cells.Item(5,18).Value = 21985.93800534225
$ws.Cells.Item(5,19).Value = 0.00571714849651024
$ws.Cells.Item(5,20).Value = 0.00571714849651024
$ws.Cells.Item(6,7).Value = 1205.102620666667
$ws.Cells.Item(6,8).Value = 3615.307862
$ws.Cells.Item(6,9).Value = 0.8667225374846176
$ws.Cells.Item(6,10).Value = 0.8667225374846176
$ws.Cells.Item(6,11).Value = 3
$ws.Cells.Item(6,12).Value = 1
$ws.Cells.Item(6,13).Value = 256.4443053333333
$ws.Cells.Item(6,14).Value = 769.332916
$ws.Cells.Item(6,15).Value = 0.8344762556643375
$ws.Cells.Item(6,16).Value = 0.8344762556643374
$ws.Cells.Item(6,17).Value = 309041.7044122429
$ws.Cells.Item(6,18).Value = 2781375.339710185
$ws.Cells.Item(6,19).Value = 0.7232593777800571
$ws.Cells.Item(6,20).Value = 0.723259377780057
$ws.Cells.Item(7,7).Value = 1205.102620666667
$ws.Cells.Item(7,8).Value = 3615.307862
$ws.Cells.Item(7,9).Value = 0.8667225374846176
$ws.Cells.Item(7,10).Value = 0.8667225374846176
$ws.Cells.Item(7,11).Value = 3
$ws.Cells.Item(7,12).Value = 1
$ws.Cells.Item(7,13).Value = 48.84026566666667
$ws.Cells.Item(7,14).Value = 146.520797
$ws.Cells.Item(7,15).Value = 0.158927459770244
$ws.Cells.Item(7,16).Value = 0.158927459770244
$ws.Cells.Item(7,17).Value = 58857.53214895624
$ws.Cells.Item(7,18).Value = 529717.7893406061
$ws.Cells.Item(7,19).Value = 0.1377460112080504
$ws.Cells.Item(7,20).Value = 0.1377460112080503
$ws.Cells.Item(8,7).Value = 180.7183073333333
$ws.Cells.Item(8,8).Value = 542.154922
$ws.Cells.Item(8,9).Value = 0.1299745160418139
$ws.Cells.Item(8,10).Value = 0.1299745160418139
$ws.Cells.Item(8,11).Value = 3
$ws.Cells.Item(8,12).Value = 1
$ws.Cells.Item(8,13).Value = 2.027115333333333
$ws.Cells.Item(8,14).Value = 6.081346
$ws.Cells.Item(8,15).Value = 0.006596284565418616
$ws.Cells.Item(8,16).Value = 0.006596284565418615
$ws.Cells.Item(8,17).Value = 366.3368518094458
$ws.Cells.Item(8,18).Value = 3297.031666285012
$ws.Cells.Item(8,19).Value = 0.0008573488940643712
$ws.Cells.Item(8,20).Value = 0.0008573488940643713
$ws.Cells.Item(9,7).Value = 180.7183073333333
$ws.Cells.Item(9,8).Value = 542.154922
$ws.Cells.Item(9,9).Value = 0.1299745160418139
$ws.Cells.Item(9,10).Value = 0.1299745160418139
$ws.Cells.Item(9,11).Value = 3
$ws.Cells.Item(9,12).Value = 1
$ws.Cells.Item(9,13).Value = 256.4443053333333
$ws.Cells.Item(9,14).Value = 769.332916
$ws.Cells.Item(9,15).Value = 0.8344762556643375
$ws.Cells.Item(9,16).Value = 0.8344762556643374
$ws.Cells.Item(9,17).Value = 46344.1807851125
$ws.Cells.Item(9,18).Value = 417097.6270660126
$ws.Cells.Item(9,19).Value = 0.1084606474783572
$ws.Cells.Item(9,20).Value = 0.1084606474783572
$ws.Cells.Item(10,7).Value = 180.7183073333333
$ws.Cells.Item(10,8).Value = 542.154922
$ws.Cells.Item(10,9).Value = 0.1299745160418139
$ws.Cells.Item(10,10).Value = 0.1299745160418139
$ws.Cells.Item(10,11).Value = 3
$ws.Cells.Item(10,12).Value = 1
$ws.Cells.Item(10,13).Value = 48.84026566666667
$ws.Cells.Item(10,14).Value = 146.520797
$ws.Cells.Item(10,15).Value = 0.158927459770244
$ws.Cells.Item(10,16).Value = 0.158927459770244
$ws.Cells.Item(10,17).Value = 8826.330140990316
$ws.Cells.Item(10,18).Value = 79436.97126891285
$ws.Cells.Item(10,19).Value = 0.02065651966939231
$ws.Cells.Item(10,20).Value = 0.02065651966939231
